$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (rows 2-7) holds date-like text; force Text format so Excel
# does not auto-convert the strings into date serial numbers.
$ws.Range("A2:A7").NumberFormat = "@"

# Header row (row 1)
$ws.Range("A1").Value = "Row"
$ws.Range("B1").Value = "Prognose"
$ws.Range("C1").Value = "surveys"
$ws.Range("D1").Value = "production"
$ws.Range("E1").Value = "orders"
$ws.Range("F1").Value = "turnover"
$ws.Range("G1").Value = "financial"
$ws.Range("H1").Value = "labor market"
$ws.Range("I1").Value = "prices"
$ws.Range("J1").Value = "national accounts"
$ws.Range("K1").Value = "Revision"

# Row 2
$ws.Range("A2").Value = "2025-09-30"
$ws.Range("B2").Value = 0.20190649201587024
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0

# Row 3
$ws.Range("A3").Value = "2025-10-15"
$ws.Range("B3").Value = 0.21305885353471507
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0.0034693103100688996
$ws.Range("E3").Value = 0.0014341519274585473
$ws.Range("F3").Value = 0.000981722311621008
$ws.Range("G3").Value = 0.000584172190987466
$ws.Range("H3").Value = -0.00007605639270772001
$ws.Range("I3").Value = -0.0010022368885525196
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0.0006089834475904998

# Row 4
$ws.Range("A4").Value = "2025-10-30"
$ws.Range("B4").Value = 0.3732122252446461
$ws.Range("C4").Value = 0.05795872952861015
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = -0.00013482287153866013
$ws.Range("F4").Value = -0.00001499200529018989
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0.0006123898564597685
$ws.Range("I4").Value = -0.011574968914945418
$ws.Range("J4").Value = 0.0027186634920311637
$ws.Range("K4").Value = -0.0016710476414749653

# Row 5
$ws.Range("A5").Value = "2025-11-15"
$ws.Range("B5").Value = 0.3825208498623729
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = -0.0001347170106651623
$ws.Range("E5").Value = -0.00007992787872255595
$ws.Range("F5").Value = 0.01113598425939579
$ws.Range("G5").Value = -0.0014779265109980898
$ws.Range("H5").Value = 0.0007335204519818629
$ws.Range("I5").Value = -0.0017763614493315888
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = -0.011951943481135763

# Row 6
$ws.Range("A6").Value = "2025-11-30"
$ws.Range("B6").Value = 0.2734205976433908
$ws.Range("C6").Value = -0.03127053134888336
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = -0.001301268262418323
$ws.Range("F6").Value = 0.0006071374760258569
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0.0011805849621470757
$ws.Range("I6").Value = -0.006992454155961061
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = -0.0007499453550176005

# Row 7
$ws.Range("A7").Value = "2025-12-15"
$ws.Range("B7").Value = 0.1990422538485066
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = -0.03178241408021746
$ws.Range("E7").Value = -0.002803394711447446
$ws.Range("F7").Value = 0.007112411248930743
$ws.Range("G7").Value = 0.002428742299472477
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0.008629798820567358

# Restore the default "Normal" style on column A so the cells do not keep
# a lingering custom number-format style index (matches original styling).
$ws.Range("A2:A7").Style = "Normal"

# Column width adjustments
$ws.Columns("D:D").ColumnWidth = 15.333333333333334
$ws.Columns("F:F").ColumnWidth = 15.0
$ws.Columns("G:G").ColumnWidth = 14.333333333333334
$ws.Columns("H:H").ColumnWidth = 15.0
$ws.Columns("I:I").ColumnWidth = 14.333333333333334
$ws.Columns("J:J").ColumnWidth = 14.166666666666666
$ws.Columns("K:K").ColumnWidth = 14.833333333333334
